# Corrige les valeurs de la colonne F (donnees non significatives) et
# met a jour la cellule active selectionnee, comme decrit dans le commit :
#  "Correction dans excel de mauvaise donnee (difference non significative)"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Nouvelles valeurs corrigees pour F2:F27 (meme ordre que les lignes 2 a 27)
$newValues = @(
    212641,
    185868,
    155865,
    129794,
    103915,
    84886,
    39657,
    39522,
    37664,
    35522,
    34923,
    35391,
    35433,
    36167,
    36947,
    37926,
    37888,
    38463,
    38797,
    39137,
    39504,
    39549,
    39889,
    40049,
    40090,
    40261
)

$arr = New-Object 'object[,]' $newValues.Length, 1
for ($i = 0; $i -lt $newValues.Length; $i++) {
    $arr[$i, 0] = $newValues[$i]
}

$ws.Range("F2:F27").Value = $arr

# La derniere cellule selectionnee dans la feuille est maintenant L20.
$ws.Range("L20").Select()
